# The deck originally had 6 slides: the real "Welcome" title slide (1),
# followed by 5 duplicate "This slide left blank for whiteboard" slides
# (2-6). The author trimmed the duplicates, keeping only a single
# whiteboard placeholder slide, so the deck ends up with 2 slides total.
$p = $ppt.ActivePresentation

# Delete slides 6, 5, 4, 3 (in descending order so earlier indices never
# shift while we still need them).
for ($i = $p.Slides.Count; $i -ge 3; $i--) {
    $p.Slides.Item($i).Delete()
}

Write-Host "Final slide count: $($p.Slides.Count)"
